# Historico.xlsx: replace the "Julian_days" index column (A) with actual
# calendar dates ("Fecha"), Jan 1 2025 (day 1) .. Aug 10 2025 (day 222).
#
# Day-of-year N (row = N + 1) maps to serial date N + 45657
# (row 2 -> 45658 == 2025-01-01, row 223 -> 45879 == 2025-08-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header. Excel will drop the now-unused "Julian_days" shared
# string and append the new "Fecha" one, shifting the other headers'
# shared-string indices down automatically.
$ws.Range("A1").Value = "Fecha"

# Replace each day-of-year value with the corresponding date serial.
for ($row = 2; $row -le 223; $row++) {
    $ws.Cells.Item($row, 1).Value = $row + 45656
}

# Apply a short-date display format to the first data cell, then fan that
# formatting out to the rest of the column via copy/paste-special so every
# cell shares the same style index (one new cellXfs entry) instead of each
# cell getting its own.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A223").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the saved selection.
[void]$ws.Range("K14").Select()
